$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing 14-pin header row (row 37) onto the
# new row so the added row matches the table's look (borders, currency
# number format on the price columns, etc.) cell-by-cell -- row 37 only
# has entries in B, F, G, H, I, K, L, M, N, so mirror exactly those.
$ws.Range("B37").Copy()
$ws.Range("B38").PasteSpecial(-4122)

$ws.Range("F37").Copy()
$ws.Range("F38").PasteSpecial(-4122)

$ws.Range("G37").Copy()
$ws.Range("G38").PasteSpecial(-4122)

$ws.Range("H37").Copy()
$ws.Range("H38").PasteSpecial(-4122)

$ws.Range("I37").Copy()
$ws.Range("I38").PasteSpecial(-4122)

$ws.Range("K37:N37").Copy()
$ws.Range("K38:N38").PasteSpecial(-4122)

# Fill in the new 24-pin 0.05" male header part data.
$ws.Range("B38").Value = "Male header"
$ws.Range("F38").Value = "CONN HEADER .050"" 24PS DL PCB AU"
$ws.Range("G38").Value = "Sullins Connector Solutions"
$ws.Range("H38").Value = "GRPB122VWVN-RC"
$ws.Range("I38").Value = "S9015E-12-ND"
$ws.Range("K38").Value = 1.93
$ws.Range("L38").Value = 1.596
$ws.Range("M38").Value = 1.19
$ws.Range("N38").Value = 0.884

# Leave the selection where the author's session ended up.
$ws.Range("N39").Select()
